$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the "Cooper" label typo (remove leading space, correct spelling) for rows 2-39
for ($r = 2; $r -le 39; $r++) {
    $ws.Cells.Item($r, 9).Value = "Copper"
}

# Update sheet view: zoom to 70% and select I2:I39 (mirrors the saved workbook view)
$ws.Select()
$excel.ActiveWindow.Zoom = 70
$ws.Range("I2:I39").Select()
